$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (appears on Overview!E2:F3 and on each locale sheet's Status column C2:C3)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Handback info for zh-cn: target file link, handback file name, handback
#    datetime.
# ---------------------------------------------------------------------------
$aUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ae9495e52c78566d65e16f09790c8c92c691dbc5/e2e/a.md"

$wsZhCn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-01 04:41:40"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $aUrl, [Type]::Missing, [Type]::Missing, "a.md")

$wsZhCn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-01 04:41:40"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $aUrl, [Type]::Missing, [Type]::Missing, "a.md")

# ---------------------------------------------------------------------------
# 3. Handback info for de-de: target file link, handback file name, handback
#    datetime.
# ---------------------------------------------------------------------------
$wsDeDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-01 04:41:48"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $aUrl, [Type]::Missing, [Type]::Missing, "a.md")

$wsDeDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-01 04:41:48"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $aUrl, [Type]::Missing, [Type]::Missing, "a.md")
